# Update "想去人数" (want-to-go count) values in column F on both the
# "展览" sheet and the consolidated "全部类型" sheet to reflect newly
# scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) -> row : new value
$sheet1Updates = @{
    2  = 6776
    4  = 429
    5  = 70
    8  = 102
    13 = 413
    15 = 1617
    17 = 3404
    19 = 229
    21 = 2036
    22 = 143
    25 = 2
    28 = 137
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheet4.xml) -> row : new value
$sheet4Updates = @{
    2  = 6776
    4  = 429
    5  = 70
    9  = 102
    14 = 413
    16 = 1617
    18 = 3404
    20 = 229
    22 = 2036
    23 = 143
    26 = 2
    29 = 137
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
